# Insert a new price-record row at row 387 (pushing existing rows 387-440
# down to 388-441), then populate the newly inserted row with the new
# "Murcott / Primera" record for Provincia de Limarí.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row before the current row 387; this shifts rows
# 387:440 down to 388:441 and extends the sheet dimension to A1:T441.
$ws.Rows.Item(387).Insert()

# Populate the new row 387 with the new record.
$ws.Cells.Item(387, 1).Value  = 9
$ws.Cells.Item(387, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(387, 3).Value  = "Metropolitana"
$ws.Cells.Item(387, 4).Value  = 44617
$ws.Cells.Item(387, 5).Value  = 13
$ws.Cells.Item(387, 6).Value  = "Fruta"
$ws.Cells.Item(387, 7).Value  = 100102
$ws.Cells.Item(387, 8).Value  = "Cítricos"
$ws.Cells.Item(387, 9).Value  = 100102004
$ws.Cells.Item(387, 10).Value = "Mandarina"
$ws.Cells.Item(387, 11).Value = "Murcott"
$ws.Cells.Item(387, 12).Value = "Primera"
$ws.Cells.Item(387, 13).Value = 140
$ws.Cells.Item(387, 14).Value = 12000
$ws.Cells.Item(387, 15).Value = 12000
$ws.Cells.Item(387, 16).Value = 12000
$ws.Cells.Item(387, 17).Value = "$/caja 15 kilos"
$ws.Cells.Item(387, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(387, 19).Value = 800
$ws.Cells.Item(387, 20).Value = 15
